# Populate the "Bleach" worksheet with the pulp-bleaching process rows
# (meta-process table: KnownQty / k_QtyFrom / UnknownQty / u_QtyTo / Calculation / Variable).
# Mirrors the layout already used on the "Pulp" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bleach")

# Row 2: unbleached pulp -> bleached pulp (product ratio)
$ws.Cells.Item(2,2).Value = "unbleached pulp"
$ws.Cells.Item(2,3).Value = "inflows"
$ws.Cells.Item(2,4).Value = "bleached pulp"
$ws.Cells.Item(2,5).Value = "outflows"
$ws.Cells.Item(2,6).Value = "ratio"
$ws.Cells.Item(2,7).Value = "pulp ratio"

# Row 3: bleached pulp -> filter cake
$ws.Cells.Item(3,2).Value = "bleached pulp"
$ws.Cells.Item(3,3).Value = "outflows"
$ws.Cells.Item(3,4).Value = "filter cake"
$ws.Cells.Item(3,5).Value = "outflows"
$ws.Cells.Item(3,6).Value = "ratio"
$ws.Cells.Item(3,7).Value = "filter cake ratio"

# Row 4: bleached pulp -> effluent
$ws.Cells.Item(4,2).Value = "bleached pulp"
$ws.Cells.Item(4,3).Value = "outflows"
$ws.Cells.Item(4,4).Value = "effluent"
$ws.Cells.Item(4,5).Value = "outflows"
$ws.Cells.Item(4,6).Value = "ratio"
$ws.Cells.Item(4,7).Value = "effluent ratio"

# Row 5: bleached pulp -> washing water demand
$ws.Cells.Item(5,2).Value = "bleached pulp"
$ws.Cells.Item(5,3).Value = "outflows"
$ws.Cells.Item(5,4).Value = "washing water"
$ws.Cells.Item(5,5).Value = "inflows"
$ws.Cells.Item(5,6).Value = "ratio"
$ws.Cells.Item(5,7).Value = "washing water demand"

# Row 6: bleached pulp -> NaOH demand
$ws.Cells.Item(6,2).Value = "bleached pulp"
$ws.Cells.Item(6,3).Value = "outflows"
$ws.Cells.Item(6,4).Value = "NaOH"
$ws.Cells.Item(6,5).Value = "inflows"
$ws.Cells.Item(6,6).Value = "ratio"
$ws.Cells.Item(6,7).Value = "NaOH demand"

# Row 7: bleached pulp -> H2O2 demand
$ws.Cells.Item(7,2).Value = "bleached pulp"
$ws.Cells.Item(7,3).Value = "outflows"
$ws.Cells.Item(7,4).Value = "H2O2"
$ws.Cells.Item(7,5).Value = "inflows"
$ws.Cells.Item(7,6).Value = "ratio"
$ws.Cells.Item(7,7).Value = "H2O2 demand"

# Row 8: bleached pulp -> O2 demand
$ws.Cells.Item(8,2).Value = "bleached pulp"
$ws.Cells.Item(8,3).Value = "outflows"
$ws.Cells.Item(8,4).Value = "O2"
$ws.Cells.Item(8,5).Value = "inflows"
$ws.Cells.Item(8,6).Value = "ratio"
$ws.Cells.Item(8,7).Value = "O2 demand"

# Row 9: bleached pulp -> MgSO4 demand
$ws.Cells.Item(9,2).Value = "bleached pulp"
$ws.Cells.Item(9,3).Value = "outflows"
$ws.Cells.Item(9,4).Value = "MgSO4"
$ws.Cells.Item(9,5).Value = "inflows"
$ws.Cells.Item(9,6).Value = "ratio"
$ws.Cells.Item(9,7).Value = "MgSO4 demand"

# Row 10: bleached pulp -> Talc demand
$ws.Cells.Item(10,2).Value = "bleached pulp"
$ws.Cells.Item(10,3).Value = "outflows"
$ws.Cells.Item(10,4).Value = "Talc"
$ws.Cells.Item(10,5).Value = "inflows"
$ws.Cells.Item(10,6).Value = "ratio"
$ws.Cells.Item(10,7).Value = "Talc demand"

# Row 11: bleached pulp -> NaCl3 demand
$ws.Cells.Item(11,2).Value = "bleached pulp"
$ws.Cells.Item(11,3).Value = "outflows"
$ws.Cells.Item(11,4).Value = "NaCl3"
$ws.Cells.Item(11,5).Value = "inflows"
$ws.Cells.Item(11,6).Value = "ratio"
$ws.Cells.Item(11,7).Value = "NaCl3 demand"

# Row 12: bleached pulp -> H2SO4 demand
$ws.Cells.Item(12,2).Value = "bleached pulp"
$ws.Cells.Item(12,3).Value = "outflows"
$ws.Cells.Item(12,4).Value = "H2SO4"
$ws.Cells.Item(12,5).Value = "inflows"
$ws.Cells.Item(12,6).Value = "ratio"
$ws.Cells.Item(12,7).Value = "H2SO4 demand"

# Row 13: bleached pulp -> Methanol demand
$ws.Cells.Item(13,2).Value = "bleached pulp"
$ws.Cells.Item(13,3).Value = "outflows"
$ws.Cells.Item(13,4).Value = "Methanol"
$ws.Cells.Item(13,5).Value = "inflows"
$ws.Cells.Item(13,6).Value = "ratio"
$ws.Cells.Item(13,7).Value = "Methanol demand"

# Row 14: bleached pulp -> electricity demand
$ws.Cells.Item(14,2).Value = "bleached pulp"
$ws.Cells.Item(14,3).Value = "outflows"
$ws.Cells.Item(14,4).Value = "electricity"
$ws.Cells.Item(14,5).Value = "inflows"
$ws.Cells.Item(14,6).Value = "ratio"
$ws.Cells.Item(14,7).Value = "electricity demand"

# Row 15: electricity -> CONSUMED electricity (returnvalue)
$ws.Cells.Item(15,2).Value = "electricity"
$ws.Cells.Item(15,3).Value = "inflows"
$ws.Cells.Item(15,4).Value = "CONSUMED electricity"
$ws.Cells.Item(15,5).Value = "outflows"
$ws.Cells.Item(15,6).Value = "returnvalue"

# The "Pulp" sheet was the previously active/selected sheet; move the
# live selection there (no longer the active tab) before activating Bleach.
$wsPulp = $wb.Worksheets.Item("Pulp")
$wsPulp.Range("B10:G11").Select()

# "Bleach" is now the active (front) worksheet with its own live selection.
$ws.Activate()
$ws.Range("G10").Select()
